$d = $word.ActiveDocument

$d.Content.Find.Execute("Resume - Patti Fernandez", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Lebenslauf: Patti Fernandez", 2)

$d.Content.Find.Execute("ABC Studios: Lead Animator (Jan 2018 - Present)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ABC Studios: Lead Animator (Jan. 2018 - heute)", 2)

$d.Content.Find.Execute("XYZ Media: Senior Animator (Jun 2015 - Dez 2017)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "XYZ Media: Senior Animator (Jun. 2015 - Dez. 2017)", 2)

$d.Content.Find.Execute("MNO Entertainment: Junior Animator (Sep 2012 - Mai 2015)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "MNO Entertainment: Junior Animator (Sept. 2012 - Mai 2015)", 2)
